$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1737  # was 1734
$ws.Cells.Item(3, 6).Value = 10133  # was 10124
$ws.Cells.Item(8, 6).Value = 1619  # was 1616
$ws.Cells.Item(9, 6).Value = 175  # was 172
$ws.Cells.Item(14, 6).Value = 477  # was 475
$ws.Cells.Item(15, 6).Value = 1176  # was 1175
$ws.Cells.Item(16, 6).Value = 128  # was 127
$ws.Cells.Item(20, 6).Value = 347  # was 346
$ws.Cells.Item(24, 6).Value = 1155  # was 1154
$ws.Cells.Item(25, 6).Value = 696  # was 695
$ws.Cells.Item(27, 6).Value = 38  # was 37
$ws.Cells.Item(29, 6).Value = 227  # was 225
$ws.Cells.Item(31, 6).Value = 421  # was 402
$ws.Cells.Item(33, 6).Value = 372  # was 371
$ws.Cells.Item(35, 6).Value = 615  # was 609
$ws.Cells.Item(36, 6).Value = 739  # was 736
$ws.Cells.Item(38, 6).Value = 1270  # was 1268
$ws.Cells.Item(41, 6).Value = 330  # was 331
$ws.Cells.Item(43, 6).Value = 354  # was 353
$ws.Cells.Item(46, 6).Value = 77  # was 76

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(16, 6).Value = 17  # was 18
$ws.Cells.Item(18, 6).Value = 1085  # was 1084
$ws.Cells.Item(20, 6).Value = 634  # was 614
$ws.Cells.Item(24, 6).Value = 73  # was 74
$ws.Cells.Item(28, 6).Value = 364  # was 363

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 2516  # was 2515
$ws.Cells.Item(7, 6).Value = 4067  # was 4064
$ws.Cells.Item(8, 6).Value = 61  # was 59
$ws.Cells.Item(10, 6).Value = 293  # was 292
$ws.Cells.Item(11, 6).Value = 185  # was 186

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1737  # was 1734
$ws.Cells.Item(4, 6).Value = 10133  # was 10124
$ws.Cells.Item(7, 6).Value = 4067  # was 4064
$ws.Cells.Item(8, 6).Value = 61  # was 59
$ws.Cells.Item(9, 6).Value = 293  # was 292
$ws.Cells.Item(10, 6).Value = 293  # was 292
$ws.Cells.Item(12, 6).Value = 1619  # was 1616
$ws.Cells.Item(13, 6).Value = 175  # was 172
$ws.Cells.Item(19, 6).Value = 1176  # was 1175
$ws.Cells.Item(20, 6).Value = 128  # was 127
$ws.Cells.Item(26, 6).Value = 1085  # was 1084
$ws.Cells.Item(27, 6).Value = 347  # was 346
$ws.Cells.Item(30, 6).Value = 1155  # was 1154
$ws.Cells.Item(31, 6).Value = 696  # was 695
$ws.Cells.Item(32, 6).Value = 73  # was 74
$ws.Cells.Item(34, 6).Value = 227  # was 225
$ws.Cells.Item(35, 6).Value = 364  # was 363
$ws.Cells.Item(36, 6).Value = 421  # was 403
$ws.Cells.Item(38, 6).Value = 372  # was 371
$ws.Cells.Item(40, 6).Value = 615  # was 609
$ws.Cells.Item(42, 6).Value = 739  # was 736
$ws.Cells.Item(46, 6).Value = 330  # was 331
$ws.Cells.Item(48, 6).Value = 354  # was 353
